$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week dates) ---
# "Volume 32   Number  35" -> "Volume 32   Number  36"
$a8 = $ws.Range("A8").Value()
$idx1 = $a8.IndexOf("35") + 1
$charsA8 = $ws.Range("A8").Characters($idx1, 2)
$charsA8.Text = "36"

# "Report Covering the Week  8/25/2025  Through  8/31/2025"
# -> "Report Covering the Week  9/1/2025  Through  9/7/2025"
$c9 = $ws.Range("C9").Value()
$idxStart = $c9.IndexOf("8/25/2025") + 1
$charsC9a = $ws.Range("C9").Characters($idxStart, 9)
$charsC9a.Text = "9/1/2025"
$c9b = $ws.Range("C9").Value()
$idxEnd = $c9b.IndexOf("8/31/2025") + 1
$charsC9b = $ws.Range("C9").Characters($idxEnd, 9)
$charsC9b.Text = "9/7/2025"

# --- Crime statistics table updates ---
$ws.Range("G14").Value = 1
$ws.Range("I14").Value = 3
$ws.Range("K14").Value = 50
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = -75
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -35.714285714285
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 8
$ws.Range("H16").Value = -20
$ws.Range("I16").Value = 68
$ws.Range("J16").Value = 84
$ws.Range("K16").Value = -19.047619047619
$ws.Range("L16").Value = -13.924050632911
$ws.Range("M16").Value = -56.129032258064
$ws.Range("N16").Value = -85.313174946004
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -66.666666666666
$ws.Range("I17").Value = 101
$ws.Range("J17").Value = 127
$ws.Range("K17").Value = -20.472440944881
$ws.Range("L17").Value = -10.619469026548
$ws.Range("M17").Value = 14.772727272727
$ws.Range("N17").Value = -56.837606837606
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -44.444444444444
$ws.Range("I18").Value = 45
$ws.Range("J18").Value = 48
$ws.Range("K18").Value = -6.25
$ws.Range("L18").Value = -31.818181818181
$ws.Range("M18").Value = -76.190476190476
$ws.Range("N18").Value = -94.346733668341
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -8.333333333333
$ws.Range("F19").Value = 42
$ws.Range("H19").Value = -22.222222222222
$ws.Range("I19").Value = 354
$ws.Range("J19").Value = 403
$ws.Range("K19").Value = -12.158808933002
$ws.Range("L19").Value = -23.043478260869
$ws.Range("M19").Value = 8.923076923076
$ws.Range("N19").Value = -22.538293216630
$ws.Range("C20").Value = 5
$ws.Range("E20").Value = 400
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 60
$ws.Range("I20").Value = 106
$ws.Range("J20").Value = 121
$ws.Range("K20").Value = -12.396694214876
$ws.Range("L20").Value = 24.705882352941
$ws.Range("M20").Value = -7.017543859649
$ws.Range("N20").Value = -94.470526864893
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -26.666666666666
$ws.Range("F21").Value = 79
$ws.Range("H21").Value = -27.522935779816
$ws.Range("I21").Value = 686
$ws.Range("J21").Value = 792
$ws.Range("K21").Value = -13.383838383838
$ws.Range("L21").Value = -15.828220858895
$ws.Range("M21").Value = -22.573363431151
$ws.Range("N21").Value = -82.396715422119
$ws.Range("C23").NumberFormat = '#,##0'
$ws.Range("C23").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 15
$ws.Range("K23").Value = -37.5
$ws.Range("L23").Value = -16.666666666666
$ws.Range("M23").Value = -31.818181818181
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = -57.5
$ws.Range("F24").Value = 103
$ws.Range("G24").Value = 141
$ws.Range("H24").Value = -26.950354609929
$ws.Range("I24").Value = 863
$ws.Range("J24").Value = 1105
$ws.Range("K24").Value = -21.900452488687
$ws.Range("L24").Value = 0.935672514619
$ws.Range("M24").Value = 27.098674521354
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 31
$ws.Range("E25").Value = -48.387096774193
$ws.Range("F25").Value = 77
$ws.Range("G25").Value = 114
$ws.Range("H25").Value = -32.456140350877
$ws.Range("I25").Value = 618
$ws.Range("J25").Value = 914
$ws.Range("K25").Value = -32.385120350109
$ws.Range("L25").Value = 1.145662847790
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 15
$ws.Range("H26").Value = 106.666666666667
$ws.Range("I26").Value = 214
$ws.Range("J26").Value = 219
$ws.Range("K26").Value = -2.283105022831
$ws.Range("L26").Value = 13.227513227513
$ws.Range("M26").Value = -13.360323886639
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 2
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 3
$ws.Range("J27").Value = 12
$ws.Range("K27").Value = -25
$ws.Range("L27").Value = -50
$ws.Range("D28").Value = 1
$ws.Range("G28").Value = 5
$ws.Range("J28").Value = 26
$ws.Range("K28").Value = -34.615384615384
$ws.Range("M29").Value = -61.538461538461
$ws.Range("M30").Value = -50
$ws.Range("D31").NumberFormat = '#,##0'
$ws.Range("D31").Value = 1
$ws.Range("E31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E31").Value = -100
$ws.Range("G31").NumberFormat = '#,##0'
$ws.Range("G31").Value = 1
$ws.Range("H31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H31").Value = -100
$ws.Range("J31").Value = 10
$ws.Range("K31").Value = -80
